# "Update - Avec les n°" - add a numbered "N°" column (G) to the sheet,
# marking reserved / removed rows, fix a misspelled professor surname.

function RGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Fix the typo "Soiukal" -> "Soukal" for the professor name in D45.
# ---------------------------------------------------------------------
$ws.Range("D45").Value = "Soukal"

# ---------------------------------------------------------------------
# New column G width (≈ default width, explicit per the target layout).
# ---------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 11.42578125

# ---------------------------------------------------------------------
# Header cell G1: "N°"  (right aligned, vertical centered, wrapped, Arial 10)
# ---------------------------------------------------------------------
$hdr = $ws.Range("G1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 10
$hdr.HorizontalAlignment = -4152   # xlHAlignRight
$hdr.VerticalAlignment = -4108     # xlVAlignCenter
$hdr.WrapText = $true
$hdr.Value = "N°"

# ---------------------------------------------------------------------
# Sequence-number cells: right aligned, vertical centered, wrapped, Arial 10,
# no fill. Applied across every contiguous block of G column that is a
# plain running number (i.e. everything except the RESERVE/SUPPRIME rows).
# ---------------------------------------------------------------------
$numRanges = @("G2:G6", "G8:G45", "G47:G51")
foreach ($addr in $numRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.HorizontalAlignment = -4152   # xlHAlignRight
    $rng.VerticalAlignment = -4108     # xlVAlignCenter
    $rng.WrapText = $true
}

$numbers = @{
    2=1; 3=2; 4=3; 5=4; 6=5;
    8=6; 9=7; 10=8; 11=9; 12=10; 13=11; 14=12; 15=13; 16=14; 17=15;
    18=16; 19=17; 20=18; 21=19; 22=20; 23=21; 24=22; 25=23; 26=24;
    27=25; 28=26; 29=27; 30=28; 31=29; 32=30; 33=31; 34=32; 35=33;
    36=34; 37=35; 38=36; 39=37; 40=38; 41=39; 42=40; 43=41; 44=42; 45=43;
    47=45; 48=46; 49=47; 50=48
}
foreach ($row in $numbers.Keys) {
    $ws.Cells.Item($row, 7).Value = $numbers[$row]
}

# ---------------------------------------------------------------------
# G7: "RESERVE" - Arial 10, vertical centered, wrapped, blue fill,
# default (general/left) horizontal alignment.
# ---------------------------------------------------------------------
$reserve = $ws.Range("G7")
$reserve.Font.Name = "Arial"
$reserve.Font.Size = 10
$reserve.VerticalAlignment = -4108   # xlVAlignCenter
$reserve.WrapText = $true
$reserve.Interior.Color = RGB 0 0 255
$reserve.Value = "RESERVE"

# ---------------------------------------------------------------------
# G46: "SUPPRIME" - Arial 10, vertical centered, wrapped, red fill,
# default (general/left) horizontal alignment.
# ---------------------------------------------------------------------
$suppr = $ws.Range("G46")
$suppr.Font.Name = "Arial"
$suppr.Font.Size = 10
$suppr.VerticalAlignment = -4108     # xlVAlignCenter
$suppr.WrapText = $true
$suppr.Interior.Color = RGB 255 0 0
$suppr.Value = "SUPPRIME"

# ---------------------------------------------------------------------
# G51: 44 - right aligned, vertical centered, wrapped, orange fill.
# ---------------------------------------------------------------------
$g51 = $ws.Range("G51")
$g51.Interior.Color = RGB 255 153 0
$g51.Value = 44

# ---------------------------------------------------------------------
# Selection cursor ends up on D45 (matches the saved view state).
# ---------------------------------------------------------------------
$ws.Range("D45").Select() | Out-Null
